# Add two new scoring rows (D/ST Points-Against and Yards-Against) to the
# ScoringTable, matching the "generated PA/YA weekly files" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lo = $ws.ListObjects.Item(1)

# Grow the table by two rows - this keeps the table ref, autoFilter ref and
# sheet dimension in sync automatically.
$row60 = $lo.ListRows.Add()
$row61 = $lo.ListRows.Add()

# Row 60: D/ST - PA (points allowed)
$ws.Range("A60").Value = "D/ST"
$ws.Range("B60").Value = "PA"
$ws.Range("C60").Value = "PA"
$ws.Range("D60").Value = 0

# Row 61: D/ST - YA (yards allowed)
$ws.Range("A61").Value = "D/ST"
$ws.Range("B61").Value = "YA"
$ws.Range("C61").Value = "YA"
$ws.Range("D61").Value = 0

# Match the author's new selection / scroll position after adding the rows.
$ws.Range("I57").Select()

# Page setup was touched (portrait orientation) in the same save.
$ws.PageSetup.Orientation = 1
